$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 137
$ws.Cells.Item(6, 6).Value = 277
$ws.Cells.Item(7, 6).Value = 13187
$ws.Cells.Item(8, 6).Value = 69
$ws.Cells.Item(10, 6).Value = 292
$ws.Cells.Item(11, 6).Value = 4506
$ws.Cells.Item(12, 6).Value = 6804
$ws.Cells.Item(15, 6).Value = 3544
$ws.Cells.Item(16, 6).Value = 44
$ws.Cells.Item(19, 6).Value = 178
$ws.Cells.Item(24, 6).Value = 3689
$ws.Cells.Item(26, 6).Value = 4000
$ws.Cells.Item(27, 6).Value = 4000
$ws.Cells.Item(29, 6).Value = 1931
$ws.Cells.Item(30, 6).Value = 109
$ws.Cells.Item(31, 6).Value = 255
$ws.Cells.Item(32, 6).Value = 6942
$ws.Cells.Item(35, 6).Value = 2045
$ws.Cells.Item(36, 6).Value = 2050
$ws.Cells.Item(37, 6).Value = 1306
$ws.Cells.Item(39, 6).Value = 1092
$ws.Cells.Item(47, 6).Value = 149
$ws.Cells.Item(48, 6).Value = 1846
$ws.Cells.Item(49, 6).Value = 76

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(15, 6).Value = 106

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 478
$ws.Cells.Item(3, 6).Value = 649

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 137
$ws.Cells.Item(5, 6).Value = 478
$ws.Cells.Item(6, 6).Value = 649
$ws.Cells.Item(8, 6).Value = 277
$ws.Cells.Item(9, 6).Value = 13187
$ws.Cells.Item(10, 6).Value = 69
$ws.Cells.Item(13, 6).Value = 292
$ws.Cells.Item(14, 6).Value = 4506
$ws.Cells.Item(15, 6).Value = 6804
$ws.Cells.Item(17, 6).Value = 3544
$ws.Cells.Item(18, 6).Value = 44
$ws.Cells.Item(28, 6).Value = 4000
$ws.Cells.Item(30, 6).Value = 109
$ws.Cells.Item(31, 6).Value = 255
$ws.Cells.Item(32, 6).Value = 6942
$ws.Cells.Item(33, 6).Value = 106
$ws.Cells.Item(36, 6).Value = 2046
$ws.Cells.Item(37, 6).Value = 2050
$ws.Cells.Item(38, 6).Value = 1306
$ws.Cells.Item(40, 6).Value = 1092
$ws.Cells.Item(45, 6).Value = 149
$ws.Cells.Item(47, 6).Value = 1847
$ws.Cells.Item(48, 6).Value = 76
